# Apply the "cryptos list" update (GitHub Actions style refresh) to Sheet1.
# Columns: A=index, B=Coin, C=Link, D=Price, E=Volume(1h)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Rows 33 & 34: "InternetComputer(DFINITY)" and "Maker" swap positions
# (Maker moves up to row 33, InternetComputer moves down to row 34),
# each also getting a freshly refreshed Price/Volume reading.
# ---------------------------------------------------------------------
$ws.Cells.Item(33, 2).Value = "Maker"
$ws.Cells.Item(33, 3).Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Cells.Item(33, 4).Value = "1.408.00"
$ws.Cells.Item(33, 5).Value = "  -1.10%  "

$ws.Cells.Item(34, 2).Value = "InternetComputer(DFINITY)"
$ws.Cells.Item(34, 3).Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "3.09"
$ws.Cells.Item(34, 5).Value = "  -0.13%  "

# ---------------------------------------------------------------------
# Price (column D) refreshes. The source data is always plain text
# (e.g. "23.42", "3.40", "0.1000"); many of these look like numbers and
# Excel would otherwise auto-convert them - silently stripping trailing
# zeros or introducing floating point noise (23.42 -> 23.420000000000002).
# Cells whose new reading contains exactly one '.' are at risk of this,
# so force those to Text format before writing so the literal string is
# preserved, matching the original inline-string content exactly.
# Readings with two '.' (thousand separator style, e.g. "27.936.51")
# are never parsed as numbers so they are left with default formatting.
# ---------------------------------------------------------------------
$priceUpdates = @{
    2  = "27.936.51"
    3  = "1.633.09"
    4  = "0.998"
    5  = "211.87"
    8  = "23.42"
    12 = "1.864.31"
    13 = "1.630.98"
    14 = "4.05"
    15 = "0.561"
    16 = "65.53"
    17 = "27.935.52"
    18 = "232.32"
    20 = "7.54"
    21 = "0.997"
    22 = "10.38"
    23 = "4.35"
    25 = "154.45"
    28 = "15.64"
    32 = "3.40"
    40 = "0.870"
    42 = "0.997"
    43 = "67.03"
    44 = "1.83"
    47 = "1.774.21"
    48 = "88.13"
    50 = "0.1000"
    51 = "0.0505"
}

foreach ($row in $priceUpdates.Keys) {
    $newValue = $priceUpdates[$row]
    $cell = $ws.Cells.Item($row, 4)
    $dotCount = $newValue.Split('.').Count - 1
    if ($dotCount -eq 1) {
        $cell.NumberFormat = "@"
    }
    $cell.Value = $newValue
}

# ---------------------------------------------------------------------
# Volume(1h) (column E) refreshes - always plain text percentages, never
# at risk of numeric auto-conversion because of the leading/trailing
# spaces and percent sign.
# ---------------------------------------------------------------------
$volumeUpdates = @{
    2  = "  +0.11%  "
    3  = "  -0.69%  "
    4  = "  -0.24%  "
    5  = "  -0.77%  "
    6  = "  -0.49%  "
    7  = "  -0.35%  "
    8  = "  -0.52%  "
    9  = "  -2.10%  "
    10 = "  -0.41%  "
    11 = "  +0.31%  "
    12 = "  -0.69%  "
    13 = "  -0.85%  "
    14 = "  -0.43%  "
    15 = "  -2.19%  "
    16 = "  -0.53%  "
    17 = "  +0.16%  "
    18 = "  +0.83%  "
    19 = "  -0.01%  "
    20 = "  -1.17%  "
    21 = "  -0.36%  "
    22 = "  -4.54%  "
    23 = "  -1.06%  "
    24 = "  -3.22%  "
    25 = "  +1.25%  "
    26 = "  +0.04%  "
    27 = "  -0.88%  "
    28 = "  -0.47%  "
    29 = "  -0.35%  "
    30 = "  -0.93%  "
    31 = "  -0.73%  "
    32 = "  +1.91%  "
    35 = "  -0.26%  "
    36 = "  +9.11%  "
    37 = "  +0.61%  "
    38 = "  +1.62%  "
    39 = "  +0.12%  "
    40 = "  -1.73%  "
    41 = "  -1.59%  "
    42 = "  -0.37%  "
    43 = "  -2.27%  "
    44 = "  +1.73%  "
    45 = "  +0.49%  "
    46 = "  -0.52%  "
    47 = "  -0.63%  "
    48 = "  -0.85%  "
    49 = "  -3.39%  "
    50 = "  -0.66%  "
    51 = "  -0.15%  "
}

foreach ($row in $volumeUpdates.Keys) {
    $ws.Cells.Item($row, 5).Value = $volumeUpdates[$row]
}
